$wb = $excel.ActiveWorkbook

# --- Sheet 1 (01_IB전략컨설팅부): remove the last two rows (퀄리타스반도체, 워트) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(29).Delete()
$ws1.Rows.Item(29).Delete()

# --- Sheet 2 (02_38커뮤니케이션(최근일자기준)): insert a new deal row at the top of the
#     data (에이피알) and drop the oldest row (에이텀) that falls off the bottom ---
$ws2 = $wb.Worksheets.Item(2)

# Insert a fresh (unformatted) row at row 3 - this shifts the old row3..row21 block
# down to row4..row22 natively (keeping their original cell types/number formats
# intact), while leaving row 2's original content untouched in place.
$ws2.Rows.Item(3).Insert()

# Move the former row 2 (케이웨더 ...) down into the row 3 slot that was just opened,
# via a real range copy so text-that-looks-numeric (e.g. "2000") keeps its string type.
$ws2.Range("A2:F2").Copy($ws2.Range("A3:F3"))

# Write the new top deal into row 2.
$ws2.Cells.Item(2, 1).Value = "에이피알"
$ws2.Cells.Item(2, 2).Value = "2024.01.22~01.26"
$ws2.Cells.Item(2, 3).Value = "147,000~200,000"
$ws2.Cells.Item(2, 4).Value = "-"
$ws2.Cells.Item(2, 5).Value = 55713
$ws2.Cells.Item(2, 6).Value = "신한투자증권,하나증권"

# The whole table grew by one row (22 rows of data + header); drop the oldest deal
# (에이텀) that rolled off the bottom of the tracked window.
$ws2.Rows.Item(22).Delete()
